$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.904.13'
$ws.Range('E2').Value = '  -1.89%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.432.16'
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.38'
$ws.Range('E5').Value = '  -2.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.63'
$ws.Range('E6').Value = '  -1.80%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.498'
$ws.Range('E8').Value = '  -2.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.430.88'
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.147'
$ws.Range('E10').Value = '  -6.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.164'
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.332'
$ws.Range('E12').Value = '  -5.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.72'
$ws.Range('E13').Value = '  -3.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.883.78'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '67.889.38'
$ws.Range('E15').Value = '  -1.92%  '
$ws.Range('E16').Value = '  -5.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '22.85'
$ws.Range('E17').Value = '  -5.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.389.07'
$ws.Range('E18').Value = '  -3.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.71'
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '337.07'
$ws.Range('E20').Value = '  -2.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.99'
$ws.Range('E21').Value = '  -4.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.70'
$ws.Range('E22').Value = '  -3.84%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -4.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.90'
$ws.Range('E25').Value = '  -5.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.59'
$ws.Range('E27').Value = '  -6.79%  '
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.97'
$ws.Range('E29').Value = '  -7.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0814'
$ws.Range('E30').Value = '  -6.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.02'
$ws.Range('E31').Value = '  -8.10%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '415.82'
$ws.Range('E33').Value = '  -5.63%  '
$ws.Range('E34').Value = '  -5.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.61'
$ws.Range('E35').Value = '  -4.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '157.22'
$ws.Range('E36').Value = '  +1.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.98'
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  -4.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.65'
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.298'
$ws.Range('E41').Value = '  -4.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.27'
$ws.Range('E42').Value = '  -6.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.45'
$ws.Range('E43').Value = '  -7.45%  '
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '132.36'
$ws.Range('E45').Value = '  -4.26%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.00'
$ws.Range('E46').Value = '  -6.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.28'
$ws.Range('E47').Value = '  -4.05%  '
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.470'
$ws.Range('E49').Value = '  -7.72%  '
$ws.Range('E50').Value = '  -3.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0901'
$ws.Range('E51').Value = '  -1.79%  '
